$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Aman Malviya"
$ws.Range("B1").Select()
